$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(5, 1).Value = "2024-10-24 14:45:28"
$ws.Cells.Item(5, 2).Value = 0.9985354969052938
$ws.Cells.Item(5, 3).Value = 0.00807231991497321
$ws.Cells.Item(5, 4).Value = 0.0001208275992999146
$ws.Cells.Item(5, 5).Value = 0.01099216081122882
$ws.Cells.Item(5, 6).Value = 0.9871104351019552
$ws.Cells.Item(5, 7).Value = 0.003673089070717696
$ws.Cells.Item(5, 8).Value = 0.00003688071937762312
$ws.Cells.Item(5, 9).Value = 0.006072949808587513
